$d = $word.ActiveDocument

$replacements = @(
    @("89×88=7832", "19×16=304"),
    @("60×85=5100", "58×20=1160"),
    @("95×57=5415", "74×70=5180"),
    @("57×41=2337", "14×75=1050"),
    @("41×53=2173", "37×65=2405"),
    @("69×37=2553", "31×14=434"),
    @("17×74=1258", "44×41=1804"),
    @("52×44=2288", "13×64=832"),
    @("63×41=2583", "44×70=3080"),
    @("67×66=4422", "43×56=2408"),
    @("16×29=464",  "74×59=4366"),
    @("96×31=2976", "19×19=361"),
    @("59×28=1652", "37×61=2257"),
    @("44×81=3564", "91×40=3640"),
    @("24×85=2040", "17×51=867"),
    @("55×11=605",  "33×54=1782"),
    @("72×74=5328", "94×67=6298"),
    @("54×44=2376", "62×31=1922"),
    @("69×52=3588", "23×66=1518"),
    @("71×34=2414", "16×53=848"),
    @("14×74=1036", "25×58=1450"),
    @("50×99=4950", "29×66=1914"),
    @("93×11=1023", "92×64=5888"),
    @("78×72=5616", "69×25=1725"),
    @("66×38=2508", "38×64=2432")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
